$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Three tables (slides 14, 15, 16) switch from table style
#    {5E579BFD-AA1A-498D-B762-5F5B03FDCECE} to
#    {53901B89-36F5-40FE-844C-04DDB4058114}.
# ---------------------------------------------------------------------------
$newStyleId = "{53901B89-36F5-40FE-844C-04DDB4058114}"
$tableSlideIndexes = @(14, 15, 16)
foreach ($idx in $tableSlideIndexes) {
    $slide = $p.Slides.Item($idx)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Theme colour scheme swap: the "Integral" (Red Violet) scheme that is
#    currently applied to the deck's master is replaced by the classic
#    "Office Theme" (Office) colour values.
#    Order exposed by ThemeColorScheme: dk1, lt1, dk2, lt2, accent1-6,
#    hlink, folHlink. Values below are the RGB() (BGR-packed) integers for
#    the target Office Theme palette.
# ---------------------------------------------------------------------------
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$theme = $p.SlideMaster.Theme
$themeColors = $theme.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i - 1]
}
